# "Memcached - Workload C test 1,2,3"
#
# Fill in the raw per-test (1,2,3) numbers for Memcached's Workload C
# (Read-only) block, which lets the existing AVERAGEA()/"=D" formulas in
# column C resolve (they were #DIV/0! before because D:F were empty).
# Also repoint the "Workloads - data" summary sheet's Memcached column
# from the old "test 1" column (D) to the new averaged column (C), and
# extend it down into the Workload B / Workload C blocks that previously
# had no Memcached entry at all.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Memcached sheet: Workload C (rows 45-53) raw test data, columns D/E/F
# ---------------------------------------------------------------------
$memcached = $wb.Worksheets.Item("Memcached")

$memcached.Range("D46").Value = 221716
$memcached.Range("E46").Value = 221312
$memcached.Range("F46").Value = 221032

$memcached.Range("D47").Value = 4510.2744050948004
$memcached.Range("E47").Value = 4518.5078079814903
$memcached.Range("F47").Value = 4524.2317854428302

$memcached.Range("D48").Value = 1000000
$memcached.Range("E48").Value = 1000000
$memcached.Range("F48").Value = 1000000

$memcached.Range("D49").Value = 881.60456099999999
$memcached.Range("E49").Value = 880.207448
$memcached.Range("F49").Value = 878.70895800000005

$memcached.Range("D50").Value = 385
$memcached.Range("E50").Value = 354
$memcached.Range("F50").Value = 389

$memcached.Range("D51").Value = 65279
$memcached.Range("E51").Value = 22719
$memcached.Range("F51").Value = 29215

$memcached.Range("D52").Value = 1235
$memcached.Range("E52").Value = 1237
$memcached.Range("F52").Value = 1228

$memcached.Range("D53").Value = 1605
$memcached.Range("E53").Value = 1629
$memcached.Range("F53").Value = 1612

# ---------------------------------------------------------------------
# 2. "Workloads - data" sheet: repoint existing Memcached column F
#    (Workload A, rows 4-17) from column D to column C.
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Workloads - data")

for ($row = 4; $row -le 17; $row++) {
    $memRow = $row + 4
    $data.Range("F$row").Formula = "=Memcached!C$memRow"
}

# ---------------------------------------------------------------------
# 3. "Workloads - data" sheet: add the missing Memcached column F for
#    the Workload B block (rows 23-36) and the Workload C block
#    (rows 42-49), mirroring the same column-C reference pattern.
# ---------------------------------------------------------------------
for ($row = 23; $row -le 36; $row++) {
    $memRow = $row + 4
    $data.Range("F$row").Formula = "=Memcached!C$memRow"
}

for ($row = 42; $row -le 49; $row++) {
    $memRow = $row + 4
    $data.Range("F$row").Formula = "=Memcached!C$memRow"
}

# ---------------------------------------------------------------------
# 4. Restore the recorded selection on each sheet (each sheet keeps its
#    own <selection> regardless of which one ends up active/on top).
# ---------------------------------------------------------------------
$redis = $wb.Worksheets.Item("Redis")
$redis.Activate()
$redis.Range("F50").Select()

$data.Activate()
$data.Range("F25").Select()

# Memcached becomes the active/visible tab last, matching activeTab="4"
# and its tabSelected="1" in the saved file.
$memcached.Activate()
$memcached.Range("C35").Select()
